# Auto-generated Excel COM-interop script to apply the crypto price update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, $NewValue)
    # Force the cell to stay a text cell (matches original inlineStr/shared-string text type)
    # instead of letting Excel auto-convert numeric-looking strings to numbers.
    $Cell.NumberFormat = "@"
    $Cell.Value = $NewValue
    $Cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "71.270.89"
Set-TextCell $ws.Range("E2") "  +6.35%  "
# Row 3
Set-TextCell $ws.Range("D3") "3.794.04"
Set-TextCell $ws.Range("E3") "  +22.57%  "
# Row 4
Set-TextCell $ws.Range("E4") "  +0.15%  "
# Row 5
Set-TextCell $ws.Range("D5") "618.85"
Set-TextCell $ws.Range("E5") "  +7.95%  "
# Row 6
Set-TextCell $ws.Range("D6") "180.83"
Set-TextCell $ws.Range("E6") "  +1.59%  "
# Row 7
Set-TextCell $ws.Range("D7") "3.791.71"
Set-TextCell $ws.Range("E7") "  +22.58%  "
# Row 8
Set-TextCell $ws.Range("E8") "  +0.08%  "
# Row 9
Set-TextCell $ws.Range("E9") "  +5.88%  "
# Row 10
Set-TextCell $ws.Range("E10") "  +9.14%  "
# Row 11
Set-TextCell $ws.Range("D11") "6.53"
Set-TextCell $ws.Range("E11") "  +2.49%  "
# Row 12
Set-TextCell $ws.Range("E12") "  +7.28%  "
# Row 13
Set-TextCell $ws.Range("D13") "40.42"
Set-TextCell $ws.Range("E13") "  +12.17%  "
# Row 14
Set-TextCell $ws.Range("E14") "  +6.85%  "
# Row 15
Set-TextCell $ws.Range("D15") "4.444.90"
Set-TextCell $ws.Range("E15") "  +23.10%  "
# Row 16
Set-TextCell $ws.Range("D16") "3.803.66"
Set-TextCell $ws.Range("E16") "  +22.84%  "
# Row 17
Set-TextCell $ws.Range("D17") "71.394.77"
Set-TextCell $ws.Range("E17") "  +6.59%  "
# Row 18
Set-TextCell $ws.Range("E18") "  +1.54%  "
# Row 19
Set-TextCell $ws.Range("E19") "  +7.85%  "
# Row 20
Set-TextCell $ws.Range("D20") "520.87"
Set-TextCell $ws.Range("E20") "  +6.91%  "
# Row 21
Set-TextCell $ws.Range("D21") "16.92"
Set-TextCell $ws.Range("E21") "  +1.14%  "
# Row 22
Set-TextCell $ws.Range("D22") "9.40"
Set-TextCell $ws.Range("E22") "  +21.65%  "
# Row 23
Set-TextCell $ws.Range("D23") "0.751"
Set-TextCell $ws.Range("E23") "  +9.53%  "
# Row 24
Set-TextCell $ws.Range("B24") "Fetch.AI"
Set-TextCell $ws.Range("C24") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws.Range("D24") "2.54"
Set-TextCell $ws.Range("E24") "  +12.03%  "
# Row 25
Set-TextCell $ws.Range("B25") "Litecoin"
Set-TextCell $ws.Range("C25") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell $ws.Range("D25") "88.92"
Set-TextCell $ws.Range("E25") "  +6.62%  "
# Row 26
Set-TextCell $ws.Range("D26") "13.57"
Set-TextCell $ws.Range("E26") "  +7.84%  "
# Row 27
Set-TextCell $ws.Range("D27") "11.13"
Set-TextCell $ws.Range("E27") "  +9.62%  "
# Row 28
Set-TextCell $ws.Range("E28") "  +0.09%  "
# Row 29
Set-TextCell $ws.Range("E29") "  +10.75%  "
# Row 30
Set-TextCell $ws.Range("D30") "8.11"
Set-TextCell $ws.Range("E30") "  +2.95%  "
# Row 31
Set-TextCell $ws.Range("D31") "2.90"
Set-TextCell $ws.Range("E31") "  +12.11%  "
# Row 32
Set-TextCell $ws.Range("B32") "EthereumClassic"
Set-TextCell $ws.Range("C32") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws.Range("D32") "32.56"
Set-TextCell $ws.Range("E32") "  +15.86%  "
# Row 33
Set-TextCell $ws.Range("B33") "PEPE"
Set-TextCell $ws.Range("C33") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws.Range("D33") "0.0000113"
Set-TextCell $ws.Range("E33") "  +19.94%  "
# Row 34
Set-TextCell $ws.Range("E34") "  +4.70%  "
# Row 35
Set-TextCell $ws.Range("E35") "  +0.13%  "
# Row 36
Set-TextCell $ws.Range("E36") "  +12.06%  "
# Row 37
Set-TextCell $ws.Range("D37") "6.17"
Set-TextCell $ws.Range("E37") "  +10.82%  "
# Row 38
Set-TextCell $ws.Range("E38") "  +10.85%  "
# Row 39
Set-TextCell $ws.Range("D39") "0.344"
Set-TextCell $ws.Range("E39") "  +9.90%  "
# Row 40
Set-TextCell $ws.Range("E40") "  +9.50%  "
# Row 41
Set-TextCell $ws.Range("D41") "51.68"
Set-TextCell $ws.Range("E41") "  +5.29%  "
# Row 42
Set-TextCell $ws.Range("D42") "442.82"
Set-TextCell $ws.Range("E42") "  +20.17%  "
# Row 43
Set-TextCell $ws.Range("D43") "3.181.13"
Set-TextCell $ws.Range("E43") "  +13.64%  "
# Row 44
Set-TextCell $ws.Range("B44") "Cosmos"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws.Range("D44") "8.87"
Set-TextCell $ws.Range("E44") "  +7.85%  "
# Row 45
Set-TextCell $ws.Range("B45") "Arweave"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextCell $ws.Range("D45") "44.45"
Set-TextCell $ws.Range("E45") "  -6.48%  "
# Row 46
Set-TextCell $ws.Range("D46") "2.82"
Set-TextCell $ws.Range("E46") "  +4.16%  "
# Row 47
Set-TextCell $ws.Range("D47") "0.0366"
Set-TextCell $ws.Range("E47") "  +6.12%  "
# Row 48
Set-TextCell $ws.Range("D48") "28.01"
Set-TextCell $ws.Range("E48") "  +9.72%  "
# Row 49
Set-TextCell $ws.Range("D49") "140.20"
Set-TextCell $ws.Range("E49") "  +3.25%  "
# Row 50
Set-TextCell $ws.Range("E50") "  +0.02%  "
# Row 51
Set-TextCell $ws.Range("E51") "  +7.90%  "
